$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values that look numeric (e.g. "1.01", "214.27") but
# must remain plain text, matching the inlineStr cells in the source workbook.
# Force a text format before assignment (cell by cell - multi-area ranges do
# not reliably propagate NumberFormat) so Excel does not coerce the strings to
# numbers, then clear the temporary formatting afterwards so each cell style
# returns to the workbook default (no explicit "s" attribute).
$priceCellAddrs = @("D2", "D3", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D20", "D21", "D22", "D23", "D25", "D26", "D28", "D29", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D50", "D51")
foreach ($addr in $priceCellAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.002.59"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.632.15"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").Value = "214.27"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").Value = "0.501"
$ws.Range("E6").Value = "  -1.30%  "
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "0.250"
$ws.Range("E8").Value = "  -2.53%  "
$ws.Range("D9").Value = "0.0618"
$ws.Range("E9").Value = "  -3.36%  "
$ws.Range("D10").Value = "18.14"
$ws.Range("E10").Value = "  -7.68%  "
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("D12").Value = "1.863.74"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").Value = "1.633.63"
$ws.Range("E13").Value = "  -3.62%  "
$ws.Range("D14").Value = "4.16"
$ws.Range("E14").Value = "  -3.08%  "
$ws.Range("D15").Value = "0.524"
$ws.Range("E15").Value = "  -3.92%  "
$ws.Range("D16").Value = "25.988.36"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").Value = "0.0₃0738"
$ws.Range("E17").Value = "  -3.86%  "
$ws.Range("D18").Value = "61.22"
$ws.Range("E18").Value = "  -3.47%  "
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Value = "190.12"
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("D21").Value = "4.22"
$ws.Range("E21").Value = "  -2.97%  "
$ws.Range("D22").Value = "9.63"
$ws.Range("E22").Value = "  -3.02%  "
$ws.Range("D23").Value = "6.06"
$ws.Range("E23").Value = "  -2.75%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "143.89"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").Value = "1.78"
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("D28").Value = "6.72"
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("D29").Value = "15.12"
$ws.Range("E29").Value = "  -2.69%  "
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("D31").Value = "0.0479"
$ws.Range("E31").Value = "  -3.73%  "
$ws.Range("D32").Value = "3.11"
$ws.Range("E32").Value = "  -4.92%  "
$ws.Range("D33").Value = "3.11"
$ws.Range("E33").Value = "  -5.45%  "
$ws.Range("E34").Value = "  -1.75%  "
$ws.Range("D35").Value = "1.47"
$ws.Range("E35").Value = "  -3.90%  "
$ws.Range("D36").Value = "1.124.31"
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("D37").Value = "0.856"
$ws.Range("E37").Value = "  -5.53%  "
$ws.Range("D38").Value = "2.43"
$ws.Range("E38").Value = "  -1.00%  "
$ws.Range("D39").Value = "0.516"
$ws.Range("E39").Value = "  -4.74%  "
$ws.Range("D40").Value = "0.0154"
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("D41").Value = "97.98"
$ws.Range("E41").Value = "  -1.60%  "
$ws.Range("D42").Value = "0.773"
$ws.Range("E42").Value = "  -2.99%  "
$ws.Range("D43").Value = "1.774.09"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").Value = "5.22"
$ws.Range("E44").Value = "  -4.97%  "
$ws.Range("D45").Value = "0.0₆0112"
$ws.Range("E45").Value = "  -4.20%  "
$ws.Range("D46").Value = "54.74"
$ws.Range("E46").Value = "  -3.36%  "
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").Value = "1.47"
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").Value = "1.01"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.47"
$ws.Range("E51").Value = "  -3.43%  "

foreach ($addr in $priceCellAddrs) {
    $ws.Range($addr).ClearFormats()
}

